# Updated grading rubrics, clarified instructions
$wb     = $excel.ActiveWorkbook
$rubric = $wb.Worksheets.Item("Rubric")

# --- Rubric sheet: insert a new blank row just above the final "Total:" row ---
# Inserting at row 23 pushes the old row 23 ("Total:" / =B11+B22) down to
# row 24, and the new blank row 23 inherits the italic sub-total formatting
# from row 22 above it.
$rubric.Rows.Item(23).Insert()

# --- Give column B some breathing room on the Rubric sheet ---
$rubric.Columns.Item(2).ColumnWidth = 6.25

# --- Rubric becomes the active / selected sheet (it was "Grade" before) ---
$rubric.Activate() | Out-Null
$rubric.Range("A1:B24").Select() | Out-Null
